# Regenerate merged AHB files
# - Rename the "_old"/"_new" header suffixes to "_FV2310"/"_FV2404"
# - Freeze the header row (top row) with a pane split
# - Wrap the data range in a native Excel Table ("Table1")

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Rename header cells: columns A-J carry the "_old" suffixed headers,
#    columns L-U carry the matching "_new" suffixed headers. Column K is
#    the untouched "diff" column.
$oldSuffixCols = @("A", "B", "C", "D", "E", "F", "G", "H", "I", "J")
$newSuffixCols = @("L", "M", "N", "O", "P", "Q", "R", "S", "T", "U")

foreach ($col in $oldSuffixCols) {
    $cell = $ws.Range($col + "1")
    $current = $cell.Value2
    $cell.Value = ($current -replace "_old$", "_FV2310")
}

foreach ($col in $newSuffixCols) {
    $cell = $ws.Range($col + "1")
    $current = $cell.Value2
    $cell.Value = ($current -replace "_new$", "_FV2404")
}

# 2) Freeze the top row (header row) in the sheet view.
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true

# 3) Turn the used range into a real Excel Table named "Table1".
$dataRange = $ws.Range("A1:U68")
$tbl = $ws.ListObjects.Add([Microsoft.Office.Interop.Excel.XlListObjectSourceType]::xlSrcRange, $dataRange, $null, [Microsoft.Office.Interop.Excel.XlYesNoGuess]::xlYes)
$tbl.Name = "Table1"
